$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "['MEC-2B-Tornearia', -, -, -]"
$ws.Range("B4").Value = "['MEC-2B-Tornearia', -, -, -]"
$ws.Range("F4").Value = "['MEC-2B-Tornearia', -, -, -]"
$ws.Range("B6").Value = "-"
$ws.Range("F6").Value = "['MEC-2B-Tornearia', -, -, -]"
$ws.Range("B7").Value = "-"
